$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("India Super League")

# Update existing row 114 odds (U114, V114)
$ws.Range("U114").Value = 1.825
$ws.Range("V114").Value = 1.975

# Copy formatting from the last existing data row (115) down onto the two
# new rows (116, 117) so that the styled cells (A: bold/border/centered,
# E: date number format) keep the same style indices used elsewhere. Only
# copy the individual styled cells (A, E) so we don't create stray blank
# cells in columns that should remain completely empty (H, I, J, AB, AC).
$ws.Range("A115").Copy()
$ws.Range("A116").PasteSpecial(-4122)
$ws.Range("A117").PasteSpecial(-4122)
$ws.Range("E115").Copy()
$ws.Range("E116").PasteSpecial(-4122)
$ws.Range("E117").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# New row 116 (id 114)
$ws.Range("A116").Value = 114
$ws.Range("B116").Value = 7749761
$ws.Range("C116").Value = "India Super League"
$ws.Range("D116").Value = "India Super League"
$ws.Range("E116").Value = 45382.45833333334
$ws.Range("F116").Value = "Mohun Bagan SG"
$ws.Range("G116").Value = "Chennaiyin FC"
$ws.Range("K116").Value = 1.533
$ws.Range("L116").Value = 4
$ws.Range("M116").Value = 5.5
$ws.Range("N116").Value = 1.533
$ws.Range("O116").Value = 4
$ws.Range("P116").Value = 5.5
$ws.Range("Q116").Value = -1
$ws.Range("R116").Value = 1.9
$ws.Range("S116").Value = 1.9
$ws.Range("T116").Value = 2.75
$ws.Range("U116").Value = 1.9
$ws.Range("V116").Value = 1.9
$ws.Range("W116").Value = 0
$ws.Range("X116").Value = 0
$ws.Range("Y116").Value = 0
$ws.Range("Z116").Value = 0
$ws.Range("AA116").Value = 0

# New row 117 (id 115)
$ws.Range("A117").Value = 115
$ws.Range("B117").Value = 7749875
$ws.Range("C117").Value = "India Super League"
$ws.Range("D117").Value = "India Super League"
$ws.Range("E117").Value = 45383.45833333334
$ws.Range("F117").Value = "Hyderabad FC"
$ws.Range("G117").Value = "Mumbai City FC"
$ws.Range("K117").Value = 9.5
$ws.Range("L117").Value = 5.5
$ws.Range("M117").Value = 1.25
$ws.Range("N117").Value = 9.5
$ws.Range("O117").Value = 5.5
$ws.Range("P117").Value = 1.25
$ws.Range("Q117").Value = 1.5
$ws.Range("R117").Value = 2
$ws.Range("S117").Value = 1.8
$ws.Range("T117").Value = 2.75
$ws.Range("U117").Value = 1.8
$ws.Range("V117").Value = 2
$ws.Range("W117").Value = 0
$ws.Range("X117").Value = 0
$ws.Range("Y117").Value = 0
$ws.Range("Z117").Value = 0
$ws.Range("AA117").Value = 0
